$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.419475
$ws.Cells.Item(2, 8).Value = 4.258425
$ws.Cells.Item(2, 9).Value = 0.1541931834006784
$ws.Cells.Item(2, 10).Value = 0.1541931834006784
$ws.Cells.Item(2, 13).Value = 16.28844733333333
$ws.Cells.Item(2, 14).Value = 48.865342
$ws.Cells.Item(2, 15).Value = 0.2176904746803693
$ws.Cells.Item(2, 16).Value = 0.2176904746803693
$ws.Cells.Item(2, 17).Value = 23.12104377848333
$ws.Cells.Item(2, 18).Value = 208.08939400635
$ws.Cells.Item(2, 19).Value = 0.03356638728697093
$ws.Cells.Item(2, 20).Value = 0.03356638728697092
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.419475
$ws.Cells.Item(3, 8).Value = 4.258425
$ws.Cells.Item(3, 9).Value = 0.1541931834006784
$ws.Cells.Item(3, 10).Value = 0.1541931834006784
$ws.Cells.Item(3, 13).Value = 27.61090666666666
$ws.Cells.Item(3, 14).Value = 82.83271999999999
$ws.Cells.Item(3, 15).Value = 0.3690119294748028
$ws.Cells.Item(3, 16).Value = 0.3690119294748029
$ws.Cells.Item(3, 17).Value = 39.19299174066666
$ws.Cells.Item(3, 18).Value = 352.736925666
$ws.Cells.Item(3, 19).Value = 0.05689912411854648
$ws.Cells.Item(3, 20).Value = 0.05689912411854648
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.419475
$ws.Cells.Item(4, 8).Value = 4.258425
$ws.Cells.Item(4, 9).Value = 0.1541931834006784
$ws.Cells.Item(4, 10).Value = 0.1541931834006784
$ws.Cells.Item(4, 13).Value = 26.266325
$ws.Cells.Item(4, 14).Value = 78.798975
$ws.Cells.Item(4, 15).Value = 0.3510419771967738
$ws.Cells.Item(4, 16).Value = 0.3510419771967739
$ws.Cells.Item(4, 17).Value = 37.284391679375
$ws.Cells.Item(4, 18).Value = 335.559525114375
$ws.Cells.Item(4, 19).Value = 0.05412827997123892
$ws.Cells.Item(4, 20).Value = 0.05412827997123892
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.419475
$ws.Cells.Item(5, 8).Value = 4.258425
$ws.Cells.Item(5, 9).Value = 0.1541931834006784
$ws.Cells.Item(5, 10).Value = 0.1541931834006784
$ws.Cells.Item(5, 13).Value = 4.658207333333333
$ws.Cells.Item(5, 14).Value = 13.974622
$ws.Cells.Item(5, 15).Value = 0.06225561864805391
$ws.Cells.Item(5, 16).Value = 0.06225561864805392
$ws.Cells.Item(5, 17).Value = 6.612208854483334
$ws.Cells.Item(5, 18).Value = 59.50987969035
$ws.Cells.Item(5, 19).Value = 0.009599392023922073
$ws.Cells.Item(5, 20).Value = 0.009599392023922073
$ws.Cells.Item(6, 9).Value = 0.3984988340349546
$ws.Cells.Item(6, 10).Value = 0.3984988340349546
$ws.Cells.Item(6, 13).Value = 16.28844733333333
$ws.Cells.Item(6, 14).Value = 48.865342
$ws.Cells.Item(6, 15).Value = 0.2176904746803693
$ws.Cells.Item(6, 16).Value = 0.2176904746803693
$ws.Cells.Item(6, 17).Value = 59.75432106784178
$ws.Cells.Item(6, 18).Value = 537.788889610576
$ws.Cells.Item(6, 19).Value = 0.08674940034064298
$ws.Cells.Item(6, 20).Value = 0.08674940034064298
$ws.Cells.Item(7, 9).Value = 0.3984988340349546
$ws.Cells.Item(7, 10).Value = 0.3984988340349546
$ws.Cells.Item(7, 13).Value = 27.61090666666666
$ws.Cells.Item(7, 14).Value = 82.83271999999999
$ws.Cells.Item(7, 15).Value = 0.3690119294748028
$ws.Cells.Item(7, 16).Value = 0.3690119294748029
$ws.Cells.Item(7, 18).Value = 911.6178192761599
$ws.Cells.Item(7, 19).Value = 0.1470508236406978
$ws.Cells.Item(7, 20).Value = 0.1470508236406978
$ws.Cells.Item(8, 9).Value = 0.3984988340349546
$ws.Cells.Item(8, 10).Value = 0.3984988340349546
$ws.Cells.Item(8, 13).Value = 26.266325
$ws.Cells.Item(8, 14).Value = 78.798975
$ws.Cells.Item(8, 15).Value = 0.3510419771967738
$ws.Cells.Item(8, 16).Value = 0.3510419771967739
$ws.Cells.Item(8, 17).Value = 96.35825841486665
$ws.Cells.Item(8, 18).Value = 867.2243257338
$ws.Cells.Item(8, 19).Value = 0.1398898186102395
$ws.Cells.Item(8, 20).Value = 0.1398898186102395
$ws.Cells.Item(9, 9).Value = 0.3984988340349546
$ws.Cells.Item(9, 10).Value = 0.3984988340349546
$ws.Cells.Item(9, 13).Value = 4.658207333333333
$ws.Cells.Item(9, 14).Value = 13.974622
$ws.Cells.Item(9, 15).Value = 0.06225561864805391
$ws.Cells.Item(9, 16).Value = 0.06225561864805392
$ws.Cells.Item(9, 17).Value = 17.08867707893511
$ws.Cells.Item(9, 18).Value = 153.798093710416
$ws.Cells.Item(9, 19).Value = 0.02480879144337426
$ws.Cells.Item(9, 20).Value = 0.02480879144337426
$ws.Cells.Item(10, 7).Value = 4.049549666666667
$ws.Cells.Item(10, 8).Value = 12.148649
$ws.Cells.Item(10, 9).Value = 0.4398900681184871
$ws.Cells.Item(10, 10).Value = 0.439890068118487
$ws.Cells.Item(10, 13).Value = 16.28844733333333
$ws.Cells.Item(10, 14).Value = 48.865342
$ws.Cells.Item(10, 15).Value = 0.2176904746803693
$ws.Cells.Item(10, 16).Value = 0.2176904746803693
$ws.Cells.Item(10, 17).Value = 65.96087646921755
$ws.Cells.Item(10, 18).Value = 593.6478882229579
$ws.Cells.Item(10, 19).Value = 0.09575987773589345
$ws.Cells.Item(10, 20).Value = 0.09575987773589344
$ws.Cells.Item(11, 7).Value = 4.049549666666667
$ws.Cells.Item(11, 8).Value = 12.148649
$ws.Cells.Item(11, 9).Value = 0.4398900681184871
$ws.Cells.Item(11, 10).Value = 0.439890068118487
$ws.Cells.Item(11, 13).Value = 27.61090666666666
$ws.Cells.Item(11, 14).Value = 82.83271999999999
$ws.Cells.Item(11, 15).Value = 0.3690119294748028
$ws.Cells.Item(11, 16).Value = 0.3690119294748029
$ws.Cells.Item(11, 17).Value = 111.8117378883644
$ws.Cells.Item(11, 18).Value = 1006.30564099528
$ws.Cells.Item(11, 19).Value = 0.1623246827932054
$ws.Cells.Item(11, 20).Value = 0.1623246827932054
$ws.Cells.Item(12, 7).Value = 4.049549666666667
$ws.Cells.Item(12, 8).Value = 12.148649
$ws.Cells.Item(12, 9).Value = 0.4398900681184871
$ws.Cells.Item(12, 10).Value = 0.439890068118487
$ws.Cells.Item(12, 13).Value = 26.266325
$ws.Cells.Item(12, 14).Value = 78.798975
$ws.Cells.Item(12, 15).Value = 0.3510419771967738
$ws.Cells.Item(12, 16).Value = 0.3510419771967739
$ws.Cells.Item(12, 17).Value = 106.3667876483083
$ws.Cells.Item(12, 18).Value = 957.3010888347749
$ws.Cells.Item(12, 19).Value = 0.1544198792615372
$ws.Cells.Item(12, 20).Value = 0.1544198792615372
$ws.Cells.Item(13, 7).Value = 4.049549666666667
$ws.Cells.Item(13, 8).Value = 12.148649
$ws.Cells.Item(13, 9).Value = 0.4398900681184871
$ws.Cells.Item(13, 10).Value = 0.439890068118487
$ws.Cells.Item(13, 13).Value = 4.658207333333333
$ws.Cells.Item(13, 14).Value = 13.974622
$ws.Cells.Item(13, 15).Value = 0.06225561864805391
$ws.Cells.Item(13, 16).Value = 0.06225561864805392
$ws.Cells.Item(13, 17).Value = 18.86364195396422
$ws.Cells.Item(13, 18).Value = 169.772777585678
$ws.Cells.Item(13, 19).Value = 0.02738562832785099
$ws.Cells.Item(13, 20).Value = 0.02738562832785099
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.068288
$ws.Cells.Item(14, 8).Value = 0.204864
$ws.Cells.Item(14, 9).Value = 0.00741791444588001
$ws.Cells.Item(14, 10).Value = 0.007417914445880009
$ws.Cells.Item(14, 13).Value = 16.28844733333333
$ws.Cells.Item(14, 14).Value = 48.865342
$ws.Cells.Item(14, 15).Value = 0.2176904746803693
$ws.Cells.Item(14, 16).Value = 0.2176904746803693
$ws.Cells.Item(14, 17).Value = 1.112305491498667
$ws.Cells.Item(14, 18).Value = 10.010749423488
$ws.Cells.Item(14, 19).Value = 0.001614809316861988
$ws.Cells.Item(14, 20).Value = 0.001614809316861988
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.068288
$ws.Cells.Item(15, 8).Value = 0.204864
$ws.Cells.Item(15, 9).Value = 0.00741791444588001
$ws.Cells.Item(15, 10).Value = 0.007417914445880009
$ws.Cells.Item(15, 13).Value = 27.61090666666666
$ws.Cells.Item(15, 14).Value = 82.83271999999999
$ws.Cells.Item(15, 15).Value = 0.3690119294748028
$ws.Cells.Item(15, 16).Value = 0.3690119294748029
$ws.Cells.Item(15, 17).Value = 1.885493594453333
$ws.Cells.Item(15, 18).Value = 16.96944235008
$ws.Cells.Item(15, 19).Value = 0.002737298922353196
$ws.Cells.Item(15, 20).Value = 0.002737298922353196
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.068288
$ws.Cells.Item(16, 8).Value = 0.204864
$ws.Cells.Item(16, 9).Value = 0.00741791444588001
$ws.Cells.Item(16, 10).Value = 0.007417914445880009
$ws.Cells.Item(16, 13).Value = 26.266325
$ws.Cells.Item(16, 14).Value = 78.798975
$ws.Cells.Item(16, 15).Value = 0.3510419771967738
$ws.Cells.Item(16, 16).Value = 0.3510419771967739
$ws.Cells.Item(16, 17).Value = 1.7936748016
$ws.Cells.Item(16, 18).Value = 16.1430732144
$ws.Cells.Item(16, 19).Value = 0.00260399935375823
$ws.Cells.Item(16, 20).Value = 0.00260399935375823
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.068288
$ws.Cells.Item(17, 8).Value = 0.204864
$ws.Cells.Item(17, 9).Value = 0.00741791444588001
$ws.Cells.Item(17, 10).Value = 0.007417914445880009
$ws.Cells.Item(17, 13).Value = 4.658207333333333
$ws.Cells.Item(17, 14).Value = 13.974622
$ws.Cells.Item(17, 15).Value = 0.06225561864805391
$ws.Cells.Item(17, 16).Value = 0.06225561864805392
$ws.Cells.Item(17, 17).Value = 0.3180996623786667
$ws.Cells.Item(17, 18).Value = 2.862896961408
$ws.Cells.Item(17, 19).Value = 0.0004618068529065961
$ws.Cells.Item(17, 20).Value = 0.0004618068529065961
